$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.366.89"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.25%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.844.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.35%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6284"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.34%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.98%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2941"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "

# Row 11
$ws.Range("E11").Value = "  -0.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.849.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.48%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.975"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.67%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6778"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.57%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001019"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.57%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "83.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.33%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.091.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.70%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.131"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "29.397.64"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.56%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "228.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.41"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.484"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.70%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1386"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.344"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.55%  "

# Row 28
$ws.Range("E28").Value = "  -0.28%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.458"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.27%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.265"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.54%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.05613"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.13%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.113"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.46%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.034"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.50%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.835"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.153"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.25%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7144"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.02%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.591"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.241.04"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.83%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01805"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "

# Row 40
$ws.Range("E40").Value = "  -0.88%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.198"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.54%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9009"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.80%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9998"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.60%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000120"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.31%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.088"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.08%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.3990"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.09%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.680"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.86%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.943"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.40%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1115"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.57%  "
